$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update text values (sharedStrings content) in place
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "NGÀNH"
$ws.Range("A2").Value = "Mã ngành"
$ws.Range("B2").Value = "Tên ngành "
$ws.Range("C2").Value = "Mã khoa"
$ws.Range("D2").Value = "Số tín chỉ của ngành"
$ws.Range("D3").Value = "creditHourTotal"
$ws.Range("D4").Value = 150

# Editing D3's text triggers an auto row-height recalculation in this runtime;
# restore the original (tiny, intentionally-hidden) row height.
$ws.Rows.Item(3).RowHeight = 0.6

# ---------------------------------------------------------------------------
# 2. Formatting tweaks
# ---------------------------------------------------------------------------
# D2 ("Số tín chỉ của ngành") adopts the same bold/red header font as A2:C2.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# D3 ("creditHourTotal") keeps bold Arial but now carries an explicit
# (theme) text colour instead of the implicit one used by A3:C3.
$ws.Range("D3").Font.ThemeColor = 1
